# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 08:05"

# Israel (row 42)
$ws.Range("B42").Value = 17008
$ws.Range("C42").Value = 21
$ws.Range("E42").Value = 1948

# Afganistan (row 48)
$ws.Range("B48").Value = 14429
$ws.Range("C48").Value = 770
$ws.Range("D48").Value = 1303
$ws.Range("E48").Value = 12878
$ws.Range("G48").Value = 2
$ws.Range("H48").Value = 248

# Uzbekistan (row 78)
$ws.Range("B78").Value = 3513
$ws.Range("C78").Value = 45
$ws.Range("E78").Value = 771

# Bulgaria (row 86)
$ws.Range("B86").Value = 2499
$ws.Range("C86").Value = 14
$ws.Range("D86").Value = 1064
$ws.Range("E86").Value = 1296
$ws.Range("G86").Value = 3
$ws.Range("H86").Value = 139

# Taiwan (row 142)
$ws.Range("D142").Value = 421
$ws.Range("E142").Value = 14

# Camboya (row 169)
$ws.Range("B169").Value = 125
$ws.Range("C169").Value = 1
$ws.Range("D169").Value = 123

# Rows 198/199: swap Fiyi <-> Curazao (with their D/H data)
$ws.Range("A198").Value = "Curazao"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 1
$ws.Range("A199").Value = "Fiyi"
$ws.Range("D199").Value = 15
$ws.Range("H199").Value = 0

# Rows 210/211: swap Seychelles <-> Montserrat (with their D/H data)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

# Rows 213/214: swap Papua Nueva Guinea <-> Islas Virgenes Britanicas (with their D/H data)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
